# Update "想去人数" (interested count) values in F column across sheets
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 122
$ws.Range("F3").Value = 1287
$ws.Range("F4").Value = 928
$ws.Range("F5").Value = 971
$ws.Range("F6").Value = 1729
$ws.Range("F7").Value = 379
$ws.Range("F8").Value = 1147
$ws.Range("F9").Value = 48
$ws.Range("F11").Value = 111
$ws.Range("F12").Value = 263
$ws.Range("F13").Value = 40
$ws.Range("F14").Value = 79
$ws.Range("F15").Value = 644
$ws.Range("F16").Value = 134
$ws.Range("F20").Value = 321
$ws.Range("F21").Value = 105
$ws.Range("F22").Value = 646
$ws.Range("F23").Value = 16
$ws.Range("F24").Value = 629
$ws.Range("F25").Value = 135
$ws.Range("F26").Value = 32
$ws.Range("F27").Value = 842
$ws.Range("F28").Value = 298
$ws.Range("F29").Value = 116
$ws.Range("F30").Value = 25
$ws.Range("F31").Value = 250
$ws.Range("F33").Value = 10
$ws.Range("F34").Value = 396

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 309
$ws.Range("F7").Value = 245
$ws.Range("F11").Value = 115
$ws.Range("F13").Value = 9

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 298

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 298
$ws.Range("F3").Value = 122
$ws.Range("F4").Value = 1287
$ws.Range("F5").Value = 928
$ws.Range("F6").Value = 971
$ws.Range("F7").Value = 1729
$ws.Range("F8").Value = 379
$ws.Range("F9").Value = 1147
$ws.Range("F10").Value = 48
$ws.Range("F13").Value = 111
$ws.Range("F14").Value = 263
$ws.Range("F15").Value = 40
$ws.Range("F16").Value = 79
$ws.Range("F17").Value = 644
$ws.Range("F18").Value = 134
$ws.Range("F22").Value = 309
$ws.Range("F25").Value = 321
$ws.Range("F27").Value = 245
$ws.Range("F28").Value = 245
$ws.Range("F29").Value = 105
$ws.Range("F30").Value = 646
$ws.Range("F31").Value = 16
$ws.Range("F32").Value = 629
$ws.Range("F33").Value = 135
$ws.Range("F34").Value = 32
$ws.Range("F35").Value = 842
$ws.Range("F36").Value = 298
$ws.Range("F39").Value = 116
$ws.Range("F40").Value = 25
$ws.Range("F41").Value = 250
$ws.Range("F43").Value = 115
$ws.Range("F44").Value = 115
$ws.Range("F46").Value = 10
$ws.Range("F48").Value = 396
$ws.Range("F49").Value = 9
